# Daily attendance processing - 2025-12-01 11:48:15
# Normalize the "Recorded By" (column G) entries so that when
# "dnasr281@gmail.com" appears together with exactly one other
# recorder (e.g. "System" or "admin@admin.com"), the other recorder
# is listed first, e.g. "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$targetUser = "dnasr281@gmail.com"
$updated = 0

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $value = $cell.Value2

    if ($value -ne $null) {
        $parts = $value -split ", "

        if ($parts.Length -eq 2 -and $parts[0] -eq $targetUser) {
            $newValue = $parts[1] + ", " + $parts[0]
            $cell.Value2 = $newValue
            $updated = $updated + 1
        }
    }
}

Write-Output "Updated $updated cells in column G"
